$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI Tnfsf13-Sdc2 LR-pair results ("Natmi following Dr Hou advice"):
# the sending/target cluster set grew from {ECs, sCs} to {ECs, FAPs, sCs}, so the
# table now has 3 x 3 = 9 data rows (rows 2-10) instead of 2 x 3 = 6 (rows 2-7),
# and every numeric column was recomputed. Columns: A Sending cluster,
# B Ligand symbol, C Receptor symbol, D Target cluster, E-T numeric stats.
$data = @(
  @("ECs","Tnfsf13","Sdc2","ECs",3,1,0.291812,0.875436,0.1179458223311005,0.1179458223311005,3,1,0.7065936666666666,2.119781,0.005187843618793344,0.005187843618793344,0.2061925110573333,1.855732599516,0.0006118844817437333,0.0006118844817437332),
  @("ECs","Tnfsf13","Sdc2","FAPs",3,1,0.291812,0.875436,0.1179458223311005,0.1179458223311005,3,1,111.9320066666667,335.79602,0.8218100075305903,0.8218100075305903,32.66310272941333,293.96792456472,0.09692905713812337,0.09692905713812336),
  @("ECs","Tnfsf13","Sdc2","sCs",3,1,0.291812,0.875436,0.1179458223311005,0.1179458223311005,3,1,23.563205,70.689615,0.1730021488506163,0.1730021488506163,6.87602597746,61.88423379714001,0.02040488071123339,0.02040488071123339),
  @("FAPs","Tnfsf13","Sdc2","ECs",2,0.6666666666666666,0.514066,1.542198,0.2077773947009016,0.2077773947009016,3,1,0.7065936666666666,2.119781,0.005187843618793344,0.005187843618793344,0.3632357798486666,3.269122018637999,0.001077916631228578,0.001077916631228578),
  @("FAPs","Tnfsf13","Sdc2","FAPs",2,0.6666666666666666,0.514066,1.542198,0.2077773947009016,0.2077773947009016,3,1,111.9320066666667,335.79602,0.8218100075305903,0.8218100075305903,57.54043893910667,517.86395045196,0.1707535423038344,0.1707535423038344),
  @("FAPs","Tnfsf13","Sdc2","sCs",2,0.6666666666666666,0.514066,1.542198,0.2077773947009016,0.2077773947009016,3,1,23.563205,70.689615,0.1730021488506163,0.1730021488506163,12.11304254153,109.01738287377,0.03594593576583863,0.03594593576583864),
  @("sCs","Tnfsf13","Sdc2","ECs",3,1,1.668241,5.004723,0.674276782967998,0.674276782967998,3,1,0.7065936666666666,2.119781,0.005187843618793344,0.005187843618793344,1.178768525073667,10.608916725663,0.003498042505821033,0.003498042505821033),
  @("sCs","Tnfsf13","Sdc2","FAPs",3,1,1.668241,5.004723,0.674276782967998,0.674276782967998,3,1,111.9320066666667,335.79602,0.8218100075305903,0.8218100075305903,186.7295627336067,1680.56606460246,0.5541274080886326,0.5541274080886326),
  @("sCs","Tnfsf13","Sdc2","sCs",3,1,1.668241,5.004723,0.674276782967998,0.674276782967998,3,1,23.563205,70.689615,0.1730021488506163,0.1730021488506163,39.309104672405,353.781942051645,0.1166513323735443,0.1166513323735443)
)

$r = 2
foreach ($row in $data) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}
